# Generate Report for handback
# Update the "Correspond Handoff Datetime" (col D) and
# "Correspond Handback DateTime" (col G) timestamps on row 4 of the
# zh-cn and de-de sheets to reflect the latest handback run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-01-18 03:34:09"
$wsZhCn.Range("G4").Value = "2016-01-18 03:35:21"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-01-18 03:34:24"
$wsDeDe.Range("G4").Value = "2016-01-18 03:35:43"
